# Limpieza y arreglo de interfaz
# Insert a new aggregated "AMBOS SEXOS" row above the existing "Hombres" row,
# pushing the existing data rows down. The new row's numeric values are the
# sum of the two existing data rows (Hombres + Mujeres). The last existing
# data row (Mujeres) is dropped, keeping the sheet at 3 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values of rows 2 and 3 before we overwrite anything.
$row2 = @{}
$row3 = @{}
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row3[$col] = $ws.Range("${col}3").Value2
}

# New row 2: sum of old row2 + old row3, label "AMBOS SEXOS"
foreach ($col in $cols) {
    $ws.Range("${col}4").Value = $ws.Range("${col}3").Value2
}
foreach ($col in $cols) {
    $ws.Range("${col}3").Value = $row2[$col]
}
foreach ($col in $cols) {
    $ws.Range("${col}2").Value = ($row2[$col] + $row3[$col])
}

$ws.Range("V2").Value = "AMBOS SEXOS"
$ws.Range("V3").Value = "AMBOS SEXOS"
$ws.Range("V4").Value = "AMBOS SEXOS"

$wb.Save()
